# unify the conception of DataNode, DataTable, Entity.
#
# - "Record_Building" sheet is removed entirely (its data now lives
#   elsewhere / is unused), which also drops its exclusively-used
#   shared strings, comments and VML drawing automatically on save.
# - Remaining sheets are renamed to the new Data* naming convention:
#     Property1                  -> DataNode_1
#     Property2                  -> DataNode_2
#     Record_Hero                -> DataTable_Hero
#     Record_Bag                 -> DataTable_Bag
#     Record_CommPropertyValue   -> DataTable_CommPropertyValue
#     Record_Task                -> DataTable_Task
#     Component                  (unchanged)
# - The active/selected tab moves from the first sheet to DataTable_Hero.

$wb = $excel.ActiveWorkbook

# Avoid the "are you sure you want to delete" prompt when removing the sheet.
$excel.DisplayAlerts = $false

# Drop the Record_Building worksheet entirely.
$wb.Worksheets("Record_Building").Delete() | Out-Null

# Rename the remaining sheets to the unified Data* naming scheme.
$wb.Worksheets("Property1").Name = "DataNode_1"
$wb.Worksheets("Property2").Name = "DataNode_2"
$wb.Worksheets("Record_Hero").Name = "DataTable_Hero"
$wb.Worksheets("Record_Bag").Name = "DataTable_Bag"
$wb.Worksheets("Record_CommPropertyValue").Name = "DataTable_CommPropertyValue"
$wb.Worksheets("Record_Task").Name = "DataTable_Task"

$excel.DisplayAlerts = $true

# Make DataTable_Hero the active tab (was sheet index 2 after Record_Building
# is removed / equivalent to Excel's 0-based activeTab="2").
$wb.Worksheets("DataTable_Hero").Select()
